$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F (想去人数 / "want to go" count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6481
$ws1.Range("F3").Value = 115
$ws1.Range("F4").Value = 7
$ws1.Range("F5").Value = 394
$ws1.Range("F6").Value = 60
$ws1.Range("F8").Value = 535
$ws1.Range("F10").Value = 80
$ws1.Range("F12").Value = 159
$ws1.Range("F14").Value = 950
$ws1.Range("F15").Value = 3186
$ws1.Range("F16").Value = 14
$ws1.Range("F17").Value = 200
$ws1.Range("F18").Value = 1856

# Sheet "全部类型" (All types) - update column F (想去人数 / "want to go" count)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6481
$ws4.Range("F3").Value = 115
$ws4.Range("F4").Value = 7
$ws4.Range("F5").Value = 394
$ws4.Range("F6").Value = 60
$ws4.Range("F9").Value = 535
$ws4.Range("F11").Value = 80
$ws4.Range("F13").Value = 159
$ws4.Range("F15").Value = 950
$ws4.Range("F16").Value = 3186
$ws4.Range("F17").Value = 14
$ws4.Range("F18").Value = 200
$ws4.Range("F19").Value = 1856
